# Serienbrief_Performer_ITM / Config.xlsx
# "Performer ITM & SOZ angepasst"
#
# Changes applied:
#  1. Assets sheet: remove the "rpa001_inProgressFolderPath" asset row
#     (row 3) - subsequent rows shift up one place.
#  2. Constants sheet: give the existing "applicationProcessNames" setting
#     (row 5) a value of "firefox, EXCEL, WINWORD".
#  3. Restore selection / active-cell bookkeeping on each sheet so the
#     workbook re-opens the same way it was left (Assets tab stays the
#     active one).

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- 1. Assets: drop the "rpa001_inProgressFolderPath" row ---------------
$wsAssets.Rows(3).Delete()

# --- 2. Constants: applicationProcessNames value --------------------------
$wsConstants.Range("B5").Value = "firefox, EXCEL, WINWORD"

# --- 3. Selection bookkeeping ---------------------------------------------
$wsSettings.Range("A47").Select()
$wsConstants.Range("B6").Select()
$wsAssets.Range("C11").Select()
$wsAssets.Activate()
